# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. This updates the DAMSLTag (column I) and
# DialogAct (column J) values for the rows whose tags changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of Excel row number -> (DAMSLTag, DialogAct)
$updates = @{
    7   = @("aa", "Agree/Accept")
    27  = @("sv", "Statement-opinion")
    36  = @("b",  "Acknowledge (Backchannel)")
    39  = @("sd", "Statement-non-opinion")
    41  = @("sd", "Statement-non-opinion")
    66  = @("sd", "Statement-non-opinion")
    67  = @("sd", "Statement-non-opinion")
    78  = @("sd", "Statement-non-opinion")
    81  = @("sd", "Statement-non-opinion")
    94  = @("ba", "Appreciation")
    99  = @("ba", "Appreciation")
    113 = @("ba", "Appreciation")
    116 = @("sd", "Statement-non-opinion")
    117 = @("b",  "Acknowledge (Backchannel)")
    131 = @("sd", "Statement-non-opinion")
    134 = @("sd", "Statement-non-opinion")
    142 = @("sd", "Statement-non-opinion")
    144 = @("sd", "Statement-non-opinion")
    148 = @("sd", "Statement-non-opinion")
    150 = @("b",  "Acknowledge (Backchannel)")
    151 = @("sd", "Statement-non-opinion")
    156 = @("aa", "Agree/Accept")
    174 = @("aa", "Agree/Accept")
    179 = @("ba", "Appreciation")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    $ws.Cells.Item($row, 9).Value = $values[0]
    $ws.Cells.Item($row, 10).Value = $values[1]
}
